# Updated cryptos list (price + 1h volume change) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.807.67"
$ws.Range("E2").Value = "  -4.82%  "
$ws.Range("D3").Value = "3.216.26"
$ws.Range("E3").Value = "  -8.18%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'595.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'151.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.60%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.207.25"
$ws.Range("E8").Value = "  -8.31%  "
$ws.Range("D9").Value = "'0.547"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.34%  "
$ws.Range("E10").Value = "  -10.74%  "
$ws.Range("D11").Value = "'6.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.62%  "
$ws.Range("D12").Value = "'0.497"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -15.21%  "
$ws.Range("D13").Value = "'39.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -15.15%  "
$ws.Range("E14").Value = "  -11.88%  "
$ws.Range("D15").Value = "3.733.53"
$ws.Range("E15").Value = "  -8.33%  "
$ws.Range("D16").Value = "66.736.29"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "3.213.84"
$ws.Range("E17").Value = "  -8.23%  "
$ws.Range("D18").Value = "'0.114"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.57%  "
$ws.Range("D19").Value = "'533.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -13.20%  "
$ws.Range("E20").Value = "  -13.61%  "
$ws.Range("D21").Value = "'15.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -14.27%  "
$ws.Range("E22").Value = "  -13.40%  "
$ws.Range("D23").Value = "'7.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.49%  "
$ws.Range("D24").Value = "'13.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.83%  "
$ws.Range("D25").Value = "'85.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.34%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'3.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -13.74%  "
$ws.Range("E28").Value = "  -14.58%  "
$ws.Range("D29").Value = "'8.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.63%  "
$ws.Range("D30").Value = "'29.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -13.74%  "
$ws.Range("D31").Value = "'2.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.31%  "
$ws.Range("E32").Value = "  -9.40%  "
$ws.Range("D33").Value = "'548.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'6.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -18.58%  "
$ws.Range("D35").Value = "'5.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -16.15%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "'53.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.60%  "
$ws.Range("D38").Value = "'0.0425"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.74%  "
$ws.Range("D39").Value = "'9.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.41%  "
$ws.Range("D40").Value = "'0.0863"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.23%  "
$ws.Range("E41").Value = "  -12.11%  "
$ws.Range("D42").Value = "2.916.56"
$ws.Range("E42").Value = "  -13.07%  "
$ws.Range("D43").Value = "'2.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -23.91%  "
$ws.Range("D45").Value = "0.0₃0588"
$ws.Range("E45").Value = "  -19.92%  "
$ws.Range("E46").Value = "  -16.58%  "
$ws.Range("D47").Value = "'26.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -16.63%  "
$ws.Range("D48").Value = "'2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -16.02%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -12.14%  "
$ws.Range("D51").Value = "'121.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.45%  "
